$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.998.11"
$ws.Range("E2").Value = "  +2.36%  "

$ws.Range("D3").Value = "2.392.11"
$ws.Range("E3").Value = "  +3.06%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "'551.96"
$ws.Range("E5").Value = "  +2.35%  "

$ws.Range("D6").Value = "'135.73"
$ws.Range("E6").Value = "  +0.93%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "'0.571"
$ws.Range("E8").Value = "  +2.07%  "

$ws.Range("E9").Value = "  +7.22%  "

$ws.Range("D10").Value = "'5.78"
$ws.Range("E10").Value = "  +6.23%  "

$ws.Range("E11").Value = "  -1.01%  "

$ws.Range("D12").Value = "'0.360"
$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").Value = "'24.43"
$ws.Range("E13").Value = "  +3.48%  "

$ws.Range("D14").Value = "2.793.75"
$ws.Range("E14").Value = "  +1.91%  "

$ws.Range("D15").Value = "58.831.94"
$ws.Range("E15").Value = "  +2.11%  "

$ws.Range("D16").Value = "'0.0000139"
$ws.Range("E16").Value = "  +4.90%  "

$ws.Range("D17").Value = "2.370.22"
$ws.Range("E17").Value = "  +1.53%  "

$ws.Range("D18").Value = "'11.26"
$ws.Range("E18").Value = "  +6.47%  "

$ws.Range("D19").Value = "'4.37"
$ws.Range("E19").Value = "  +3.28%  "

$ws.Range("D20").Value = "'334.77"
$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").Value = "'7.01"
$ws.Range("E21").Value = "  +5.64%  "

$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").Value = "'64.65"
$ws.Range("E23").Value = "  +4.22%  "

$ws.Range("D24").Value = "'0.170"
$ws.Range("E24").Value = "  +1.39%  "

$ws.Range("D25").Value = "'0.995"
$ws.Range("E25").Value = "  -0.53%  "

$ws.Range("D26").Value = "'8.41"
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").Value = "'1.35"
$ws.Range("E27").Value = "  -3.62%  "

$ws.Range("D28").Value = "'1.79"
$ws.Range("E28").Value = "  +2.03%  "

$ws.Range("D29").Value = "0.0₃0761"
$ws.Range("E29").Value = "  +4.90%  "

$ws.Range("D30").Value = "'170.74"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("D31").Value = "'6.21"
$ws.Range("E31").Value = "  +2.11%  "

$ws.Range("D32").Value = "'18.61"
$ws.Range("E32").Value = "  +1.40%  "

$ws.Range("D33").Value = "'1.02"
$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.25"
$ws.Range("E35").Value = "  +2.45%  "

$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "'1.26"
$ws.Range("E37").Value = "  +1.56%  "

$ws.Range("D38").Value = "'40.42"
$ws.Range("E38").Value = "  +3.47%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.63"
$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.423"
$ws.Range("E40").Value = "  +12.81%  "

$ws.Range("D41").Value = "'3.74"
$ws.Range("E41").Value = "  +3.64%  "

$ws.Range("D42").Value = "'295.89"
$ws.Range("E42").Value = "  +4.31%  "

$ws.Range("D43").Value = "'141.73"
$ws.Range("E43").Value = "  -1.46%  "

$ws.Range("E44").Value = "  +2.96%  "

$ws.Range("D45").Value = "'0.0520"
$ws.Range("E45").Value = "  +4.01%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'19.06"
$ws.Range("E46").Value = "  +0.36%  "

$ws.Range("D47").Value = "'0.570"
$ws.Range("E47").Value = "  +2.25%  "

$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").Value = "'0.400"
$ws.Range("E48").Value = "  +3.95%  "

$ws.Range("D49").Value = "'0.0225"
$ws.Range("E49").Value = "  +4.62%  "

$ws.Range("D50").Value = "'11.03"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'1.57"
$ws.Range("E51").Value = "  +3.40%  "
